# Weekly fruit/vegetable price data refresh for Arveja Verde (Vega Monumental Concepcion).
# The rows get reshuffled/updated with a new week of data (old row 2 record moves to
# the bottom as row 7, a brand-new row is inserted as row 3, and all other rows shift
# up to make room), so we rewrite every data cell explicitly to the final target state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44335
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112022
$ws.Range("G2").Value = "Arveja Verde"
$ws.Range("H2").Value = "Perfection"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 32000
$ws.Range("M2").Value = 31000
$ws.Range("N2").Value = "$/malla 25 kilos"
$ws.Range("O2").Value = "Provincia de Huasco"
$ws.Range("P2").Value = 1240
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44454
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 36000
$ws.Range("L3").Value = 38000
$ws.Range("M3").Value = 37000
$ws.Range("N3").Value = "$/malla 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 1480
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44342
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100112022
$ws.Range("G4").Value = "Arveja Verde"
$ws.Range("H4").Value = "Perfection"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 32000
$ws.Range("M4").Value = 31000
$ws.Range("N4").Value = "$/malla 25 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 1240
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44399
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112022
$ws.Range("G5").Value = "Arveja Verde"
$ws.Range("H5").Value = "Perfection"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 39000
$ws.Range("L5").Value = 40000
$ws.Range("M5").Value = 39600
$ws.Range("N5").Value = "$/malla 25 kilos"
$ws.Range("O5").Value = "Provincia de Huasco"
$ws.Range("P5").Value = 1584
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44328
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = "Arveja Verde"
$ws.Range("H6").Value = "Perfection"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 33000
$ws.Range("L6").Value = 34000
$ws.Range("M6").Value = 33500
$ws.Range("N6").Value = "$/malla 25 kilos"
$ws.Range("O6").Value = "Provincia de Huasco"
$ws.Range("P6").Value = 1340
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44162
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112022
$ws.Range("G7").Value = "Arveja Verde"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("N7").Value = "$/saco 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 700
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"

# Row 7 is brand new; give its Fecha (date) cell the same date/time number format
# used by the other rows in column D (s="2" -> numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat
